# Apply the weekly update: insert two new price records (rows 49-50) for
# "Choclo" - "Choclero" at "Vega Monumental Concepción", dated 2023-02-22
# (serial 44979), pushing all subsequent rows down by two and extending the
# used range from A1:R149 to A1:R151.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 49; everything that used to
# be row 49 onward now starts at row 51.
$ws.Range("A49:A50").EntireRow.Insert()

# --- New row 49 ---------------------------------------------------------
$ws.Cells.Item(49, 1).Value = 11
$ws.Cells.Item(49, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(49, 3).Value = "Bíobío"
$ws.Cells.Item(49, 4).Value = 44979
$ws.Cells.Item(49, 5).Value = 8
$ws.Cells.Item(49, 6).Value = 100112024
$ws.Cells.Item(49, 7).Value = "Choclo"
$ws.Cells.Item(49, 8).Value = "Choclero"
$ws.Cells.Item(49, 9).Value = "Primera"
$ws.Cells.Item(49, 10).Value = 5000
$ws.Cells.Item(49, 11).Value = 450
$ws.Cells.Item(49, 12).Value = 450
$ws.Cells.Item(49, 13).Value = 450
$ws.Cells.Item(49, 14).Value = "`$/unidad"
$ws.Cells.Item(49, 15).Value = "Región Metropolitana"
$ws.Cells.Item(49, 16).Value = 450
$ws.Cells.Item(49, 17).Value = 1
$ws.Cells.Item(49, 18).Value = "Hortaliza"

# --- New row 50 ---------------------------------------------------------
$ws.Cells.Item(50, 1).Value = 11
$ws.Cells.Item(50, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(50, 3).Value = "Bíobío"
$ws.Cells.Item(50, 4).Value = 44979
$ws.Cells.Item(50, 5).Value = 8
$ws.Cells.Item(50, 6).Value = 100112024
$ws.Cells.Item(50, 7).Value = "Choclo"
$ws.Cells.Item(50, 8).Value = "Choclero"
$ws.Cells.Item(50, 9).Value = "Segunda"
$ws.Cells.Item(50, 10).Value = 2000
$ws.Cells.Item(50, 11).Value = 350
$ws.Cells.Item(50, 12).Value = 350
$ws.Cells.Item(50, 13).Value = 350
$ws.Cells.Item(50, 14).Value = "`$/unidad"
$ws.Cells.Item(50, 15).Value = "Región Metropolitana"
$ws.Cells.Item(50, 16).Value = 350
$ws.Cells.Item(50, 17).Value = 1
$ws.Cells.Item(50, 18).Value = "Hortaliza"
